$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trip report rows (13-15), appended after the existing 12 rows.
$newRows = @(
    @{ Row = 13; A = 10; B = 45693; C = 45693; D = 1; E = "대구시"; F = "대구시"; G = "참외 작목현황";                           H = "김상걸, 차수호"; I = "참외 ";    J = 45694; K = "김상걸" },
    @{ Row = 14; A = 11; B = 45693; C = 45693; D = 1; E = "경남";   F = "밀양";   G = "하우스감사 시세 동향 및 출하 상담"; H = "김용보, 이용수"; I = "감자";    J = 45694; K = "김용보" },
    @{ Row = 15; A = 12; B = 45693; C = 45693; D = 1; E = "전남";   F = "해남군"; G = "봄동배추, 대파 줄하독려";              H = "김언중 김기영";  I = "배추,대파"; J = 45694; K = "김언중" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I

    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 10).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 11).Value = $r.K
}
